# Apply progress updates to sheet "БИВТ-22-17" (first worksheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Turn several "pass" (text) marks into numeric scores of 5,
# and fill in a couple of previously empty score cells with 5.
$ws.Range("C3").Value  = 5
$ws.Range("C4").Value  = 5
$ws.Range("D11").Value = 5
$ws.Range("C16").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("C22").Value = 5
$ws.Range("F23").Value = 5
$ws.Range("D24").Value = 5
$ws.Range("D25").Value = 5

# Move the active selection to match the author's last cursor position.
$ws.Range("F25").Select()
